# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 00:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1093198
$ws.Range("C4").Value = 29004
$ws.Range("D4").Value = 151784
$ws.Range("E4").Value = 877624
$ws.Range("G4").Value = 2135
$ws.Range("H4").Value = 63790

# Row 9 - Alemania
$ws.Range("B9").Value = 163009
$ws.Range("C9").Value = 1470
$ws.Range("E9").Value = 32886
$ws.Range("G9").Value = 156
$ws.Range("H9").Value = 6623

# Row 91 - Tunez
$ws.Range("B91").Value = 994
$ws.Range("C91").Value = 14
$ws.Range("D91").Value = 305
$ws.Range("E91").Value = 648
$ws.Range("F91").Value = 24
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 41

# Row 147 - Cabo Verde
$ws.Range("B147").Value = 121
$ws.Range("C147").Value = 8
$ws.Range("D147").Value = 4
$ws.Range("E147").Value = 116
